# Revert the "update_goods_infor" merge: drop the second sheet
# (DanhSachTinhNang) that was added, rename the remaining sheet back to
# "Sheet1", and restore its prior view/selection state.

$wb = $excel.ActiveWorkbook

# Remove the "DanhSachTinhNang" worksheet entirely.
$ws2 = $wb.Worksheets.Item("DanhSachTinhNang")
$ws2.Delete()

# The only remaining sheet ("SoSanh") becomes "Sheet1".
$ws1 = $wb.Worksheets.Item("SoSanh")
$ws1.Name = "Sheet1"

# The title row (A1:E1) loses its bold styling.
$ws1.Range("A1:E1").Font.Bold = $false

# Restore the prior selection/active cell on the now-only sheet.
$ws1.Activate()
$ws1.Range("E9").Select()
